# Update the daily vaccination report workbook from 2021-01-02 to 2021-01-03.

$wb = $excel.ActiveWorkbook

# --- Rename the data sheet (also updates the defined name "Bundeslaender001"
#     that points at '01.01.21'!$A$1:$G$17 -> '02.01.21'!$A$1:$G$17) ---
$wsData = $wb.Worksheets.Item(2)
$wsData.Name = "02.01.21"

# --- Data sheet values (column layout: A=Bundesland, B=Gesamt,
#     C=Differenz zum Vortag, D=Indikation nach Alter,
#     E=Berufliche Indikation, F=Medizinische Indikation,
#     G=Pflegeheim-bewohnerIn, H=footnote markers) ---

# Row 2 - Bayern
$wsData.Range("B2").Value = 24063
$wsData.Range("C2").Value = 4014
$wsData.Range("D2").Value = 11041
$wsData.Range("E2").Value = 7428
$wsData.Range("F2").Value = 1382
$wsData.Range("G2").Value = 4200

# Row 3 - Baden-Wuerttemberg
$wsData.Range("B3").Value = 57833
$wsData.Range("C3").Value = 11742
$wsData.Range("D3").Value = 12855
$wsData.Range("E3").Value = 25636
$wsData.Range("F3").Value = 1091
$wsData.Range("G3").Value = 20428
$wsData.Range("H3").ClearContents()

# Row 4 - Berlin
$wsData.Range("B4").Value = 14616
$wsData.Range("C4").Value = 1443
$wsData.Range("D4").Value = 9817
$wsData.Range("E4").Value = 3647
$wsData.Range("G4").Value = 10968

# Row 5 - Brandenburg (footnote moves from index 33 to 34 - handled below)

# Row 6 - Bremen
$wsData.Range("B6").Value = 1837
$wsData.Range("C6").Value = 96
$wsData.Range("D6").Value = 824
$wsData.Range("E6").Value = 465
$wsData.Range("F6").Value = 9
$wsData.Range("G6").Value = 265

# Row 7 - Hamburg
$wsData.Range("B7").Value = 3704
$wsData.Range("C7").Value = 662
$wsData.Range("D7").Value = 1439
$wsData.Range("E7").Value = 2010
$wsData.Range("G7").Value = 1696

# Row 8 - Hessen
$wsData.Range("B8").Value = 30085
$wsData.Range("C8").Value = 5294
$wsData.Range("D8").Value = 8443
$wsData.Range("E8").Value = 17635
$wsData.Range("F8").Value = 1096
$wsData.Range("G8").Value = 11569

# Row 10 - Niedersachsen
$wsData.Range("B10").Value = 4962
$wsData.Range("C10").Value = 1017
$wsData.Range("D10").Value = 912
$wsData.Range("E10").Value = 2641
$wsData.Range("F10").Value = 982
$wsData.Range("G10").Value = 2918

# Row 11 - Nordrhein-Westfalen
$wsData.Range("B11").Value = 48691
$wsData.Range("C11").Value = 13142
$wsData.Range("E11").Value = 18806
$wsData.Range("F11").ClearContents()
$wsData.Range("G11").Value = 29893

# Row 12 - Rheinland-Pfalz (gets a new footnote marker, index 33 - handled below)
$wsData.Range("B12").Value = 7248
$wsData.Range("C12").Value = 59
$wsData.Range("E12").Value = 3516

# Row 13 - Saarland
$wsData.Range("B13").Value = 4149
$wsData.Range("C13").Value = 833
$wsData.Range("D13").Value = 3122
$wsData.Range("E13").Value = 495
$wsData.Range("G13").Value = 1580

# Row 14 - Sachsen
$wsData.Range("B14").Value = 4343
$wsData.Range("C14").Value = 343
$wsData.Range("D14").Value = 368
$wsData.Range("E14").Value = 3535
$wsData.Range("G14").Value = 807

# Row 15 - Sachsen-Anhalt
$wsData.Range("B15").Value = 12822
$wsData.Range("C15").Value = 1051
$wsData.Range("D15").Value = 4237
$wsData.Range("E15").Value = 6194
$wsData.Range("F15").Value = 546
$wsData.Range("G15").Value = 6410

# Row 16 - Schleswig-Holstein
$wsData.Range("B16").Value = 8933
$wsData.Range("C16").Value = 969
$wsData.Range("D16").Value = 2525
$wsData.Range("E16").Value = 5136
$wsData.Range("F16").Value = 2077
$wsData.Range("G16").Value = 3987

# --- Footnotes: H3's note is dropped entirely; H5 gets the new
#     "keine Meldung" note; H12 is newly annotated with the
#     "Nachmeldungen der mobilen Teams" note. Writing H12 before H5
#     reproduces the shared-string slot order of the target file. ---
$wsData.Range("H12").Value = "(Nachmeldungen der mobilen Teams stehen noch aus)"
$wsData.Range("H5").Value = "(keine Meldung für den 02.01.; Nachmeldung erfolgt am 04.01.2020) "

# --- View/selection state: move the tab selection + cell selection from
#     the data sheet to "Erläuterung", and update the data sheet's own
#     selection. ---
$wsData.Select()
$wsData.Range("L24:L25").Select()

$wsIntro = $wb.Worksheets.Item(1)
$wsIntro.Select()
$wsIntro.Range("E26").Select()
